# Adds "scaled" (normalized-to-first-row) columns for density, viscosity,
# thermal conductivity, kinematic viscosity and diffusivity to the Sheet1
# properties table, and a couple of extra average/format tweaks in row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Insert the four new columns (order matters: left-to-right on the
#        *current* column letters at each step) -------------------------
[void]$ws.Columns("D").Insert()   # new col D: "Scaled density"   (after Density=C)
[void]$ws.Columns("L").Insert()   # new col L: "Scaled viscosity" (after Viscosity=K)
[void]$ws.Columns("N").Insert()   # new col N: "Scaled conductivity" (after Therm. Cond.=M)
[void]$ws.Columns("P").Insert()   # new col P: "Scaled kinematic viscosity" (after Kinem. Visc.=O)

# After the inserts the table looks like:
#  A Temperature | B Pressure | C Density | D (new) | E Volume | F Internal Energy
#  | G Enthalpy | H Entropy | I Cp | J Sound Spd. | K Viscosity | L (new)
#  | M Therm. Cond. | N (new) | O Kinem. Visc. | P (new) | Q diffusivity

# --- 2. Headers (order chosen to match the shared-string table order of the
#        authored workbook: density, kin.visc., diffusivity, conductivity,
#        viscosity) ---------------------------------------------------------
$ws.Range("D1").Value = "Scaled density"
$ws.Range("P1").Value = "Scaled kinematic viscosity"
$ws.Range("R1").Value = "Scaled diffusivity"
$ws.Range("N1").Value = "Scaled conductivity"
$ws.Range("L1").Value = "Scaled viscosity"

# --- 3. Per-row formulas for the five "scaled" columns (rows 2-12) --------
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 4).Formula  = "=C$r/`$C`$2"             # D: Scaled density
    $ws.Cells.Item($r, 12).Formula = "=K$r/`$K`$2"             # L: Scaled viscosity
    $ws.Cells.Item($r, 14).Formula = "=M$r/`$M`$2"             # N: Scaled conductivity
    $ws.Cells.Item($r, 16).Formula = "=O$r/`$O`$2"             # P: Scaled kinematic viscosity
    $ws.Cells.Item($r, 18).Formula = "=Q$r/`$Q`$2"             # R: Scaled diffusivity
}

# --- 4. Row 13 additions ---------------------------------------------------
$ws.Range("C13").Formula = "=AVERAGE(C2:C12)"
$ws.Range("M13").Formula = "=AVERAGE(M2:M12)"

# O13 / Q13 already hold the (shifted) averages for kinematic viscosity and
# diffusivity; give them their own distinct number formats (Q13 first, so the
# numFmt ids are allocated in the same order as the authored workbook: 164 =
# scientific, 165 = fixed-decimal).
$ws.Range("Q13").NumberFormat = "0.000E+00"
$ws.Range("O13").NumberFormat = "0.0000000"

# P13 stays empty but picks up the "scientific" style used by its column.
$ws.Range("P13").NumberFormat = "0.00E+00"

# --- 5. Workbook-level defined name ----------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "fluid" -or $n.Name -eq "Sheet1!fluid") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$M`$12"
    }
}

# --- 6. Selection matches the authored file (M13 active) -------------------
[void]$ws.Range("M13").Select()

Write-Host "done"
